$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 5 data rows (A:E, rows 2-6) are being re-ordered into a new sequence.
# New row order (1-indexed, referring to the ORIGINAL row numbers 2..6):
#   new row2 <- old row6
#   new row3 <- old row4
#   new row4 <- old row3
#   new row5 <- old row2
#   new row6 <- old row5

$titles = @(
    "Jury finds N.B. shooter Matthew Raymond not criminally responsible for four killings",
    "P.E.I. residents in Fredericton react to shooting in New Brunswick capital",
    "Yarmouth boy's Hallelujah cover for Fredericton shooting victims is viral hit",
    "Accused in Fredericton shooting described as ‘polite and pleasant’, motive still unclear",
    "‘I saw 3 bodies laying there’: Witnesses describe chaos in fatal Fredericton shooting"
)

$timestamps = @(
    "2020-11-20T06:11:00UTC",
    "2018-08-10T15:41:01UTC",
    "2018-08-13T14:28:43UTC",
    "2018-08-12T07:41:50UTC",
    "2018-08-10T11:39:14UTC"
)

$distances = @(833, 0, 3, 2, 0)

$buckets = @("day_31_beyond", "day_0", "day_2_to_30", "day_2_to_30", "day_0")

$uris = @(
    "https://atlantic.ctvnews.ca/jury-finds-n-b-shooter-matthew-raymond-not-criminally-responsible-for-four-killings-1.5197334",
    "http://www.journalpioneer.com/news/local/pei-residents-in-fredericton-react-to-shooting-in-new-brunswick-capital-232944/",
    "https://www.cbc.ca/news/canada/nova-scotia/fredericton-shooting-joshua-cochrane-hallelujah-tribute-1.4783017",
    "https://globalnews.ca/news/4383367/fredericton-shooting-accused-motive/",
    "https://globalnews.ca/news/4380632/fredericton-shooting-police-witnesses/"
)

# Clear the existing hyperlinks on E2:E6 before rewriting values/hyperlinks,
# so Excel does not keep stale hyperlink-to-range associations around.
$ws.Range("E2:E6").Hyperlinks.Delete()

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $titles[$i]
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
    $ws.Cells.Item($row, 3).Value = $distances[$i]
    $ws.Cells.Item($row, 4).Value = $buckets[$i]
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $uris[$i])
    $ws.Cells.Item($row, 5).Value = $uris[$i]
    $ws.Cells.Item($row, 5).Style = "Hyperlink"
}
